$d = $word.ActiveDocument

function Replace-Text($old, $new) {
    $d.Content.Find.Execute($old, $true, $false, $false, $false, $false, $true, 1, $false, $new, 2)
}

# Note: order matters because "14÷3=4, 2" is both a target (old) value and
# a result (new) value of a later replacement. We must replace the cell
# that currently holds "14÷3=4, 2" first, before creating a new
# "14÷3=4, 2" elsewhere.
Replace-Text "14÷3=4, 2" "63÷5=12, 3"

Replace-Text "44÷7=6, 2" "99÷9=11, 0"
Replace-Text "70÷7=10, 0" "96÷4=24, 0"
Replace-Text "45÷4=11, 1" "52÷8=6, 4"
Replace-Text "47÷3=15, 2" "91÷3=30, 1"
Replace-Text "66÷6=11, 0" "74÷3=24, 2"
Replace-Text "89÷4=22, 1" "84÷7=12, 0"
Replace-Text "42÷5=8, 2" "69÷7=9, 6"
Replace-Text "35÷3=11, 2" "41÷2=20, 1"
Replace-Text "98÷9=10, 8" "20÷5=4, 0"
Replace-Text "65÷4=16, 1" "79÷2=39, 1"
Replace-Text "23÷4=5, 3" "36÷2=18, 0"
Replace-Text "48÷4=12, 0" "92÷2=46, 0"
Replace-Text "53÷9=5, 8" "38÷3=12, 2"
Replace-Text "76÷9=8, 4" "37÷5=7, 2"
Replace-Text "73÷3=24, 1" "58÷2=29, 0"
Replace-Text "36÷8=4, 4" "68÷4=17, 0"
Replace-Text "71÷5=14, 1" "33÷6=5, 3"
Replace-Text "92÷6=15, 2" "76÷3=25, 1"
Replace-Text "50÷5=10, 0" "97÷9=10, 7"
Replace-Text "57÷9=6, 3" "54÷3=18, 0"
Replace-Text "35÷7=5, 0" "90÷4=22, 2"
Replace-Text "11÷3=3, 2" "40÷9=4, 4"
Replace-Text "35÷9=3, 8" "14÷3=4, 2"
Replace-Text "77÷6=12, 5" "52÷3=17, 1"
